$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1113.1818
$ws.Range("J4").Value = 412.5
$ws.Range("L4").Value = 412.5
$ws.Range("N4").Value = -640.5
$ws.Range("H28").Value = 805.1667
$ws.Range("I28").Value = 915
$ws.Range("K28").Value = 915
$ws.Range("M28").Value = -430
$ws.Range("H33").Value = 28571736
$ws.Range("I33").Value = 34483108
$ws.Range("J33").Value = 97.5
$ws.Range("K33").Value = 34483108
$ws.Range("L33").Value = 97.5
$ws.Range("M33").Value = -34482879
$ws.Range("N33").Value = -555.5
$ws.Range("H64").Value = 3588.4285
$ws.Range("I64").Value = 3557.1333
$ws.Range("J64").Value = 3666.6667
$ws.Range("K64").Value = 3557.1333
$ws.Range("L64").Value = 3666.6667
$ws.Range("M64").Value = -3309.1333
$ws.Range("N64").Value = -4162.6667
$ws.Range("H67").Value = 3588.4285
$ws.Range("I67").Value = 3557.1333
$ws.Range("J67").Value = 3666.6667
$ws.Range("K67").Value = 3557.1333
$ws.Range("L67").Value = 3666.6667
$ws.Range("M67").Value = -2699.1333
$ws.Range("N67").Value = -5382.6667
$ws.Range("H98").Value = 1164.2084
$ws.Range("I98").Value = 1134.3889
$ws.Range("K98").Value = 1134.3889
$ws.Range("M98").Value = 363.6111000000001
$ws.Range("H118").Value = 2144.862
$ws.Range("J118").Value = 3007.1177
$ws.Range("L118").Value = 9021.3531
$ws.Range("N118").Value = -12335.3531
$ws.Range("H122").Value = 1164.2084
$ws.Range("I122").Value = 1134.3889
$ws.Range("K122").Value = 3403.1667
$ws.Range("M122").Value = -953.1666999999998
$ws.Range("H132").Value = 4267.6
$ws.Range("I132").Value = 1794.2413
$ws.Range("J132").Value = 10788.272
$ws.Range("K132").Value = 5382.7239
$ws.Range("L132").Value = 32364.816
$ws.Range("M132").Value = -2852.7239
$ws.Range("N132").Value = -37424.81600000001
$ws.Range("H137").Value = 5411627.5
$ws.Range("I137").Value = 5525.5
$ws.Range("J137").Value = 11771747
$ws.Range("K137").Value = 16576.5
$ws.Range("L137").Value = 35315241
$ws.Range("M137").Value = -14026.5
$ws.Range("N137").Value = -35320341
$ws.Range("H138").Value = 5557263.5
$ws.Range("I138").Value = 1232
$ws.Range("J138").Value = 11907014
$ws.Range("K138").Value = 3696
$ws.Range("L138").Value = 35721042
$ws.Range("M138").Value = 1444
$ws.Range("N138").Value = -35731322

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 16
$ws.Range("H74").Value = 13890464
$ws.Range("I74").Value = 17858258
$ws.Range("J74").Value = 3183.75
$ws.Range("K74").Value = 17858258
$ws.Range("L74").Value = 3183.75
$ws.Range("M74").Value = -17857384
$ws.Range("N74").Value = -4931.75
$ws.Range("H77").Value = 13890464
$ws.Range("I77").Value = 17858258
$ws.Range("J77").Value = 3183.75
$ws.Range("K77").Value = 89291290
$ws.Range("L77").Value = 15918.75
$ws.Range("M77").Value = -89286922
$ws.Range("N77").Value = -24654.75
$ws.Range("H132").Value = 9262365
$ws.Range("I132").Value = 13891087
$ws.Range("J132").Value = 4921.3335
$ws.Range("K132").Value = 41673261
$ws.Range("L132").Value = 14764.0005
$ws.Range("M132").Value = -41670731
$ws.Range("N132").Value = -19824.0005

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 117
$ws.Range("I7").Value = 134
$ws.Range("K7").Value = 134
$ws.Range("M7").Value = -21
$ws.Range("H31").Value = 9528830
$ws.Range("I31").Value = 4615.375
$ws.Range("K31").Value = 4615.375
$ws.Range("M31").Value = -4320.375
$ws.Range("H34").Value = 9528830
$ws.Range("I34").Value = 4615.375
$ws.Range("K34").Value = 4615.375
$ws.Range("M34").Value = -4413.375
$ws.Range("H41").Value = 19341.25
$ws.Range("I41").Value = 2300
$ws.Range("J41").Value = 25021.666
$ws.Range("K41").Value = 2300
$ws.Range("L41").Value = 25021.666
$ws.Range("M41").Value = -1872
$ws.Range("N41").Value = -25877.666
$ws.Range("H94").Value = 3672.6316
$ws.Range("I94").Value = 1298.3334
$ws.Range("J94").Value = 7742.857
$ws.Range("K94").Value = 1298.3334
$ws.Range("L94").Value = 7742.857
$ws.Range("M94").Value = -847.3334
$ws.Range("N94").Value = -8644.857
$ws.Range("H103").Value = 20291.166
$ws.Range("I103").Value = 18549.4
$ws.Range("J103").Value = 29000
$ws.Range("K103").Value = 18549.4
$ws.Range("L103").Value = 29000
$ws.Range("M103").Value = -17377.4
$ws.Range("N103").Value = -31344

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2581.818
$ws.Range("I5").Value = 733.3333
$ws.Range("K5").Value = 2199.9999
$ws.Range("M5").Value = -2087.9999
$ws.Range("H122").Value = 649.4545000000001
$ws.Range("I122").Value = 512.8570999999999
$ws.Range("J122").Value = 888.5
$ws.Range("K122").Value = 4615.7139
$ws.Range("L122").Value = 7996.5
$ws.Range("M122").Value = -2165.7139
$ws.Range("N122").Value = -12896.5
$ws.Range("H134").Value = 4087.8948
$ws.Range("I134").Value = 2747.5
$ws.Range("J134").Value = 6385.7144
$ws.Range("K134").Value = 8242.5
$ws.Range("L134").Value = 19157.1432
$ws.Range("M134").Value = -3172.5
$ws.Range("N134").Value = -29297.1432
$ws.Range("H135").Value = 2581.818
$ws.Range("I135").Value = 733.3333
$ws.Range("K135").Value = 6599.9997
$ws.Range("M135").Value = -4064.9997

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1593.75
$ws.Range("I113").Value = 1402.75
$ws.Range("J113").Value = 1784.75
$ws.Range("K113").Value = 1402.75
$ws.Range("L113").Value = 1784.75
$ws.Range("M113").Value = 767.25
$ws.Range("N113").Value = -6124.75

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1197
$ws.Range("I46").Value = 842
$ws.Range("J46").Value = 1552
$ws.Range("K46").Value = 842
$ws.Range("L46").Value = 1552
$ws.Range("M46").Value = -654
$ws.Range("N46").Value = -1928
$ws.Range("H68").Value = 2245
$ws.Range("J68").Value = 2365.7144
$ws.Range("L68").Value = 2365.7144
$ws.Range("N68").Value = -3863.7144
$ws.Range("H70").Value = 50000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 50000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 50000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -50540
$ws.Range("H71").Value = 2245
$ws.Range("J71").Value = 2365.7144
$ws.Range("L71").Value = 11828.572
$ws.Range("N71").Value = -19316.572
$ws.Range("H73").Value = 50000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 50000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 50000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -51872
$ws.Range("H122").Value = 8606.267
$ws.Range("I122").Value = 17834.666
$ws.Range("J122").Value = 6299.1665
$ws.Range("K122").Value = 53503.99800000001
$ws.Range("L122").Value = 18897.4995
$ws.Range("M122").Value = -51053.99800000001
$ws.Range("N122").Value = -23797.4995

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10131.842
$ws.Range("I62").Value = 6500.2856
$ws.Range("J62").Value = 12250.25
$ws.Range("K62").Value = 6500.2856
$ws.Range("L62").Value = 12250.25
$ws.Range("M62").Value = -5876.2856
$ws.Range("N62").Value = -13498.25
$ws.Range("H65").Value = 10131.842
$ws.Range("I65").Value = 6500.2856
$ws.Range("J65").Value = 12250.25
$ws.Range("K65").Value = 32501.428
$ws.Range("L65").Value = 61251.25
$ws.Range("M65").Value = -29381.428
$ws.Range("N65").Value = -67491.25

Write-Output "Applied all Ultima_Profits market-data updates"